# Appends 5 new match rows (rows 63-67) to the Armenia Premier League 2023-2024
# sheet, mirroring the formatting of the last existing data row (row 62):
#   - column A keeps the bold/bordered/centered "Indice" style
#   - column E keeps the custom date-time number format
# Updates the sheet dimension implicitly by writing into the new cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 63 (Indice 62): Van vs Ararat-Armenia ---
$ws.Range("A62:V62").Copy()
$ws.Range("A63").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A63").Value = 62
$ws.Range("B63").Value = 'armenia'
$ws.Range("C63").Value = 'premier-league'
$ws.Range("D63").Value = '2023-2024'
$ws.Range("E63").Value = 45227.54166666666
$ws.Range("F63").Value = 'Van'
$ws.Range("G63").Value = 0
$ws.Range("H63").Value = 'Ararat-Armenia'
$ws.Range("I63").Value = 2
$ws.Range("J63").Value = 8
$ws.Range("K63").Value = '27/10/2023 01:12'
$ws.Range("L63").Value = 10.13
$ws.Range("M63").Value = '28/10/2023 12:26'
$ws.Range("N63").Value = 5.29
$ws.Range("O63").Value = '27/10/2023 01:12'
$ws.Range("P63").Value = 5.87
$ws.Range("Q63").Value = '28/10/2023 12:26'
$ws.Range("R63").Value = 1.27
$ws.Range("S63").Value = '27/10/2023 01:12'
$ws.Range("T63").Value = 1.26
$ws.Range("U63").Value = '28/10/2023 12:26'
$ws.Range("V63").Value = 'https://www.betexplorer.com/football/armenia/premier-league/van-ararat-armenia/tUv5jeSg/'

# --- Row 64 (Indice 63): West Armenia vs Pyunik Yerevan ---
$ws.Range("A63:V63").Copy()
$ws.Range("A64").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A64").Value = 63
$ws.Range("B64").Value = 'armenia'
$ws.Range("C64").Value = 'premier-league'
$ws.Range("D64").Value = '2023-2024'
$ws.Range("E64").Value = 45228.45833333334
$ws.Range("F64").Value = 'West Armenia'
$ws.Range("G64").Value = 2
$ws.Range("H64").Value = 'Pyunik Yerevan'
$ws.Range("I64").Value = 3
$ws.Range("J64").Value = 30.52
$ws.Range("K64").Value = '29/10/2023 06:14'
$ws.Range("L64").Value = 30.52
$ws.Range("M64").Value = '29/10/2023 06:14'
$ws.Range("N64").Value = 16.62
$ws.Range("O64").Value = '29/10/2023 06:14'
$ws.Range("P64").Value = 16.62
$ws.Range("Q64").Value = '29/10/2023 06:14'
$ws.Range("R64").Value = 1.03
$ws.Range("S64").Value = '29/10/2023 06:14'
$ws.Range("T64").Value = 1.03
$ws.Range("U64").Value = '29/10/2023 06:14'
$ws.Range("V64").Value = 'https://www.betexplorer.com/football/armenia/premier-league/west-armenia-pyunik-yerevan/fXr9kFs0/'

# --- Row 65 (Indice 64): BKMA vs Ararat Yerevan ---
$ws.Range("A64:V64").Copy()
$ws.Range("A65").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A65").Value = 64
$ws.Range("B65").Value = 'armenia'
$ws.Range("C65").Value = 'premier-league'
$ws.Range("D65").Value = '2023-2024'
$ws.Range("E65").Value = 45228.58333333334
$ws.Range("F65").Value = 'BKMA'
$ws.Range("G65").Value = 2
$ws.Range("H65").Value = 'Ararat Yerevan'
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 2.81
$ws.Range("K65").Value = '28/10/2023 03:12'
$ws.Range("L65").Value = 3.23
$ws.Range("M65").Value = '29/10/2023 13:52'
$ws.Range("N65").Value = 3.07
$ws.Range("O65").Value = '28/10/2023 03:12'
$ws.Range("P65").Value = 3.22
$ws.Range("Q65").Value = '29/10/2023 13:52'
$ws.Range("R65").Value = 2.39
$ws.Range("S65").Value = '28/10/2023 03:12'
$ws.Range("T65").Value = 2.29
$ws.Range("U65").Value = '29/10/2023 13:52'
$ws.Range("V65").Value = 'https://www.betexplorer.com/football/armenia/premier-league/bkma-ararat-yerevan/n7u1iyDm/'

# --- Row 66 (Indice 65): Noah vs Alashkert ---
$ws.Range("A65:V65").Copy()
$ws.Range("A66").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A66").Value = 65
$ws.Range("B66").Value = 'armenia'
$ws.Range("C66").Value = 'premier-league'
$ws.Range("D66").Value = '2023-2024'
$ws.Range("E66").Value = 45229.58333333334
$ws.Range("F66").Value = 'Noah'
$ws.Range("G66").Value = 4
$ws.Range("H66").Value = 'Alashkert'
$ws.Range("I66").Value = 2
$ws.Range("J66").Value = 2.11
$ws.Range("K66").Value = '29/10/2023 02:12'
$ws.Range("L66").Value = 2.13
$ws.Range("M66").Value = '30/10/2023 13:59'
$ws.Range("N66").Value = 3.42
$ws.Range("O66").Value = '29/10/2023 02:12'
$ws.Range("P66").Value = 3.53
$ws.Range("Q66").Value = '30/10/2023 13:59'
$ws.Range("R66").Value = 3.08
$ws.Range("S66").Value = '29/10/2023 02:12'
$ws.Range("T66").Value = 3.29
$ws.Range("U66").Value = '30/10/2023 13:59'
$ws.Range("V66").Value = 'https://www.betexplorer.com/football/armenia/premier-league/noah-alashkert/8GtchHcs/'

# --- Row 67 (Indice 66): Urartu vs Shirak Gyumri ---
$ws.Range("A66:V66").Copy()
$ws.Range("A67").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A67").Value = 66
$ws.Range("B67").Value = 'armenia'
$ws.Range("C67").Value = 'premier-league'
$ws.Range("D67").Value = '2023-2024'
$ws.Range("E67").Value = 45229.66666666666
$ws.Range("F67").Value = 'Urartu'
$ws.Range("G67").Value = 2
$ws.Range("H67").Value = 'Shirak Gyumri'
$ws.Range("I67").Value = 2
$ws.Range("J67").Value = 1.34
$ws.Range("K67").Value = '29/10/2023 04:12'
$ws.Range("L67").Value = 1.3
$ws.Range("M67").Value = '30/10/2023 15:43'
$ws.Range("N67").Value = 4.72
$ws.Range("O67").Value = '29/10/2023 04:12'
$ws.Range("P67").Value = 5.37
$ws.Range("Q67").Value = '30/10/2023 15:59'
$ws.Range("R67").Value = 6.99
$ws.Range("S67").Value = '29/10/2023 04:12'
$ws.Range("T67").Value = 9.91
$ws.Range("U67").Value = '30/10/2023 15:59'
$ws.Range("V67").Value = 'https://www.betexplorer.com/football/armenia/premier-league/urartu-shirak-gyumri/h6NS9qTE/'

